$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (force text for plain-number-looking values, then reset style)
$ws.Range("D2").Value = "69.282.57"
$ws.Range("D3").Value = "3.683.80"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "679.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.437"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000231"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "4.304.93"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.681.41"
$ws.Range("D16").Value = "69.272.32"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.652"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "3.831.69"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "3.672.15"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0905"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "171.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.942"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000274"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.78"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) column
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E10").Value = "  -4.19%  "
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  -3.59%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -7.02%  "
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -4.81%  "
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("E43").Value = "  +4.47%  "
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("E46").Value = "  -5.69%  "
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("E48").Value = "  -4.81%  "
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("E51").Value = "  -3.99%  "
